$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.462.00'
$ws.Range('E2').Value = '  -2.67%  '
$ws.Range('D3').Value = '1.804.09'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.86%  '
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4533'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3651'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07110'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8749'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07733'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('D13').Value = '1.844.87'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.260'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.342'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '86.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.64%  '
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008561'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.009'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').Value = '26.509.07'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.969'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.38'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.978'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.988'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '112.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.831'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08652'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.042'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7272'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.423'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.110'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.92%  '
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.555'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.076'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01925'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05085'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.881'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.919'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4974'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1565'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.102'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.07%  '
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4591'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.977'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.584'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05990'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.79%  '
